$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E (PDF) - status changes from "Pendente" to "OK" for rows 2-5
$ws.Range("E2:E5").Value = "OK"

# Column F (Tempo Decorrido) - elapsed time values updated per row
$ws.Range("F2").Value = "0 hours, 0 minutes, and 15.39 seconds"
$ws.Range("F3").Value = "0 hours, 0 minutes, and 0.11 seconds"
$ws.Range("F4").Value = "0 hours, 0 minutes, and 0.10 seconds"
$ws.Range("F5").Value = "0 hours, 0 minutes, and 31.16 seconds"
